$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 891; this shifts every existing row
# (891..1001) down by one (to 892..1002), automatically growing the
# sheet from A1:R1001 to A1:R1002.
$ws.Rows("891:891").Insert()

# Populate the newly inserted row 891 with the new weekly price record.
$ws.Range("A891").Value = 8
$ws.Range("B891").Value = "Terminal La Palmera de La Serena"
$ws.Range("C891").Value = "Coquimbo"
$ws.Range("D891").Value = 45212
$ws.Range("E891").Value = 4
$ws.Range("F891").Value = 100112045
$ws.Range("G891").Value = "Zapallo"
$ws.Range("H891").Value = "Camote"
$ws.Range("I891").Value = "1a nueva(o)"
$ws.Range("J891").Value = 1500
$ws.Range("K891").Value = 1150
$ws.Range("L891").Value = 1200
$ws.Range("M891").Value = 1175
$ws.Range("N891").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O891").Value = "Perú"
$ws.Range("P891").Value = 1175
$ws.Range("Q891").Value = 1
$ws.Range("R891").Value = "Hortaliza"
